$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a brand-new "Player Info" worksheet as the very first sheet.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item(1).Activate()
$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"

$battingSheet = $wb.Worksheets.Item("ODI Batting")

# Copy the header formatting (bold, border, centered) from the existing
# "ODI Batting" header row so the new sheet's header matches the rest of
# the workbook's look & feel.
$battingSheet.Range("A1").Copy()
$playerInfo.Range("A1:D1").PasteSpecial(-4122)
$playerInfo.Application.CutCopyMode = 0

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Force the data row to be stored as text (the ID looks numeric, but the
# rest of the workbook keeps this kind of value as text too).
$playerInfo.Range("A2:D2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "4933"
$playerInfo.Range("B2").Value = "Arshdeep Singh"
$playerInfo.Range("C2").Value = "Left Handed"
$playerInfo.Range("D2").Value = "Left Arm Medium Fast"

# ---------------------------------------------------------------------------
# 2. "ODI Batting" sheet: rename MATCH_CARD_LINK -> MATCH_CODE and replace
#    the full scorecard URLs with the bare numeric match code.
# ---------------------------------------------------------------------------
$battingSheet.Range("D1").Value = "MATCH_CODE"

$battingSheet.Range("D2:D4").NumberFormat = "@"
$battingSheet.Range("D2").Value = "4669"
$battingSheet.Range("D3").Value = "4673"
$battingSheet.Range("D4").Value = "4676"

# ---------------------------------------------------------------------------
# 3. "ODI Bowling" sheet: same MATCH_CARD_LINK -> MATCH_CODE treatment.
# ---------------------------------------------------------------------------
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$bowlingSheet.Range("B1").Value = "MATCH_CODE"

$bowlingSheet.Range("B2:B3").NumberFormat = "@"
$bowlingSheet.Range("B2").Value = "4669"
$bowlingSheet.Range("B3").Value = "4676"
